$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the hyperlinks that currently live on the row-2 data cells so we
#    can rebuild them (in the new column order) with fresh relationship ids.
# ---------------------------------------------------------------------------
$ws.Range("P2:AD2").Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 2. Clear the three old-layout cells that disappear entirely in the new
#    layout (R2/S2 held the old thumbnail+manifest urls, U2 held the old
#    sort-title value - all three get relocated/renamed elsewhere).
# ---------------------------------------------------------------------------
$ws.Range("R2").Clear()
$ws.Range("S2").Clear()
$ws.Range("U2").Clear()

# ---------------------------------------------------------------------------
# 3. Re-label the header row (R1:Z1) with the new "variety of format" column
#    names, in their new order.
# ---------------------------------------------------------------------------
$ws.Range("R1").Value = "年"
$ws.Range("S1").Value = "viewingDirection"
$ws.Range("T1").Value = "ID"
$ws.Range("U1").Value = "ソート用項目"
$ws.Range("V1").Value = "機械可読ドキュメント"
$ws.Range("W1").Value = "ウェブサイトURL"
$ws.Range("X1").Value = "IIIFマニフェストURI"
$ws.Range("Y1").Value = "コレクション"
$ws.Range("Z1").Value = "サムネイル"

# ---------------------------------------------------------------------------
# 4. Write the row-2 data values into their new homes.
# ---------------------------------------------------------------------------
$ws.Range("T2").Value = "102a847a-ea96-464a-84ba-cb714696bfbd"
$ws.Range("V2").Value = "https://iiif.dl.itc.u-tokyo.ac.jp/repo/api/items/10"
$ws.Range("W2").Value = "https://iiif.dl.itc.u-tokyo.ac.jp/repo/s/koshu/"
$ws.Range("X2").Value = "https://iiif.dl.itc.u-tokyo.ac.jp/koshu/manifests/koshu.json"
$ws.Range("Y2").Value = "甲州法度之次第"
$ws.Range("Z2").Value = "https://iiif.dl.itc.u-tokyo.ac.jp/repo/files/square/10/default.jpg"

# ---------------------------------------------------------------------------
# 5. Re-create the hyperlinks, in left-to-right order, so relationship ids
#    are allocated rId1..rId10 matching the new column layout.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("P2"), "https://www.lib.j.u-tokyo.ac.jp/digitalarchive/terms/legalhistorysectionterms.html")
$ws.Hyperlinks.Add($ws.Range("Q2"), "https://iiif.dl.itc.u-tokyo.ac.jp/repo/s/law/document/102a847a-ea96-464a-84ba-cb714696bfbd")
$ws.Hyperlinks.Add($ws.Range("V2"), "https://iiif.dl.itc.u-tokyo.ac.jp/repo/api/items/10")
$ws.Hyperlinks.Add($ws.Range("W2"), "https://iiif.dl.itc.u-tokyo.ac.jp/repo/s/koshu/")
$ws.Hyperlinks.Add($ws.Range("X2"), "https://iiif.dl.itc.u-tokyo.ac.jp/koshu/manifests/koshu.json")
$ws.Hyperlinks.Add($ws.Range("Z2"), "https://iiif.dl.itc.u-tokyo.ac.jp/repo/files/square/10/default.jpg")
$ws.Hyperlinks.Add($ws.Range("AA2"), "https://iiif.dl.itc.u-tokyo.ac.jp/omekac/oa/collections/10/manifest.json")
$ws.Hyperlinks.Add($ws.Range("AB2"), "http://tapasproject.org/tapas-commons/files/甲州法度之次第")
$ws.Hyperlinks.Add($ws.Range("AC2"), "https://utda.github.io/text/rtf/102a847a-ea96-464a-84ba-cb714696bfbd.rtf")
$ws.Hyperlinks.Add($ws.Range("AD2"), "https://iiif.dl.itc.u-tokyo.ac.jp/api/iiif-search/kPzFpI4mtex7HdRmrZL1ew9r7OCgdDPvNX2g0njpVtBywlyhooghaIDF9TDMS%EF%BC%8BvUVN4VOwDy1Vp1%EF%BC%8BtH%EF%BC%8B3oU5hRhxL%EF%BC%8BEHPis3o5UucFHDPKE%3D")

# ---------------------------------------------------------------------------
# 6. Re-apply the workbook's existing hyperlink font (underline + blue) so
#    the cells land back on the original custom style instead of Excel's
#    auto-generated built-in "Hyperlink" style.
# ---------------------------------------------------------------------------
$linkCells = @("P2","Q2","V2","W2","X2","Z2","AA2","AB2","AC2","AD2")
foreach ($addr in $linkCells) {
    $ws.Range($addr).Font.Underline = $true
    $ws.Range($addr).Font.Color = 16711680
}
